{"js": "// Leave 3/9/2023 12:08 AM\n// Update the terminal-leave computation's mail-merge result text:\n//   MR -> MS\n//   INOCENCIO M.  ANGCAYA -> AMELITA V.  FERMA\n//   Casual Employee -> Agriculturist B\n//   Vice Mayor's Office Detailed At Civil Security Unit -> City Agriculture Office\n//   11,374.00 -> 26,497.00 (appears twice)\n//   92.416 -> 453.500 (appears twice)\n//   50,657.25 -> 579,102.25\nconst body = context.document.body;\n\nasync function replaceAll(findText, replaceText, matchWholeWord) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: !!matchWholeWord });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, \"Replace\");\n  }\n  await context.sync();\n}\n\nawait replaceAll(\"MR\", \"MS\", true);\nawait replaceAll(\"INOCENCIO M.  ANGCAYA\", \"AMELITA V.  FERMA\");\nawait replaceAll(\"Casual Employee\", \"Agriculturist B\");\nawait replaceAll(\"Vice Mayor's Office Detailed At Civil Security Unit\", \"City Agriculture Office\");\nawait replaceAll(\"11,374.00\", \"26,497.00\");\nawait replaceAll(\"92.416\", \"453.500\");\nawait replaceAll(\"50,657.25\", \"579,102.25\");\n", "ps1": "# Leave 3/9/2023 12:08 AM\n# Update the terminal-leave computation's mail-merge result text:\n#   MR -> MS\n#   INOCENCIO M.  ANGCAYA -> AMELITA V.  FERMA\n#   Casual Employee -> Agriculturist B\n#   Vice Mayor's Office Detailed At Civil Security Unit -> City Agriculture Office\n#   11,374.00 -> 26,497.00 (appears twice)\n#   92.416 -> 453.500 (appears twice)\n#   50,657.25 -> 579,102.25\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText, $matchWholeWord) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $findText\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $matchWholeWord\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $replaceText\n  # wdFindContinue=1, Replace:=wdReplaceAll(2)\n  $find.Execute($findText, $true, $matchWholeWord, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nReplace-Text \"MR\" \"MS\" $true\nReplace-Text \"INOCENCIO M.  ANGCAYA\" \"AMELITA V.  FERMA\" $false\nReplace-Text \"Casual Employee\" \"Agriculturist B\" $false\nReplace-Text \"Vice Mayor's Office Detailed At Civil Security Unit\" \"City Agriculture Office\" $false\nReplace-Text \"11,374.00\" \"26,497.00\" $false\nReplace-Text \"92.416\" \"453.500\" $false\nReplace-Text \"50,657.25\" \"579,102.25\" $false\n"}
